$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "With cache"
$ws.Range("B3").Value = 124.40633939999999
$ws.Range("C3").Value = 93.208470000000005

$ws.Range("A3").Select()
